$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20

$data[0,0] = "ECs"
$data[0,1] = "Wnt4"
$data[0,2] = "Fzd2"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 3.285322
$data[0,7] = 9.855966
$data[0,8] = 0.4533344065718998
$data[0,9] = 0.4533344065718997
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 0.08106100000000001
$data[0,13] = 0.243183
$data[0,14] = 0.004404799763312406
$data[0,15] = 0.004404799763312406
$data[0,16] = 0.2663114866420001
$data[0,17] = 2.396803379778
$data[0,18] = 0.001996847286769274
$data[0,19] = 0.001996847286769274

$data[1,0] = "ECs"
$data[1,1] = "Wnt4"
$data[1,2] = "Fzd2"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 3.285322
$data[1,7] = 9.855966
$data[1,8] = 0.4533344065718998
$data[1,9] = 0.4533344065718997
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 15.70489833333333
$data[1,13] = 47.114695
$data[1,14] = 0.8533935241547975
$data[1,15] = 0.8533935241547975
$data[1,16] = 51.59564800226334
$data[1,17] = 464.36083202037
$data[1,18] = 0.3868726468450174
$data[1,19] = 0.3868726468450173

$data[2,0] = "ECs"
$data[2,1] = "Wnt4"
$data[2,2] = "Fzd2"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 3.285322
$data[2,7] = 9.855966
$data[2,8] = 0.4533344065718998
$data[2,9] = 0.4533344065718997
$data[2,10] = 2
$data[2,11] = 0.6666666666666666
$data[2,12] = 0.03299466666666667
$data[2,13] = 0.098984
$data[2,14] = 0.001792907809228914
$data[2,15] = 0.001792907809228914
$data[2,16] = 0.1083981042826667
$data[2,17] = 0.975582938544
$data[2,18] = 0.0008127867977349147
$data[2,19] = 0.0008127867977349147

$data[3,0] = "ECs"
$data[3,1] = "Wnt4"
$data[3,2] = "Fzd2"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 3.285322
$data[3,7] = 9.855966
$data[3,8] = 0.4533344065718998
$data[3,9] = 0.4533344065718997
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 2.583925666666666
$data[3,13] = 7.751776999999999
$data[3,14] = 0.140408768272661
$data[3,15] = 0.1404087682726611
$data[3,16] = 8.489027839064665
$data[3,17] = 76.401250551582
$data[3,18] = 0.06365212564237818
$data[3,19] = 0.0636521256423782

$data[4,0] = "FAPs"
$data[4,1] = "Wnt4"
$data[4,2] = "Fzd2"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 1.533623
$data[4,7] = 4.600869
$data[4,8] = 0.2116212878402838
$data[4,9] = 0.2116212878402837
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 0.08106100000000001
$data[4,13] = 0.243183
$data[4,14] = 0.004404799763312406
$data[4,15] = 0.004404799763312406
$data[4,16] = 0.124317014003
$data[4,17] = 1.118853126027
$data[4,18] = 0.0009321493985907483
$data[4,19] = 0.000932149398590748

$data[5,0] = "FAPs"
$data[5,1] = "Wnt4"
$data[5,2] = "Fzd2"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 1.533623
$data[5,7] = 4.600869
$data[5,8] = 0.2116212878402838
$data[5,9] = 0.2116212878402837
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 15.70489833333333
$data[5,13] = 47.114695
$data[5,14] = 0.8533935241547975
$data[5,15] = 0.8533935241547975
$data[5,16] = 24.08539329666167
$data[5,17] = 216.768539669955
$data[5,18] = 0.1805962366161966
$data[5,19] = 0.1805962366161965

$data[6,0] = "FAPs"
$data[6,1] = "Wnt4"
$data[6,2] = "Fzd2"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 1.533623
$data[6,7] = 4.600869
$data[6,8] = 0.2116212878402838
$data[6,9] = 0.2116212878402837
$data[6,10] = 2
$data[6,11] = 0.6666666666666666
$data[6,12] = 0.03299466666666667
$data[6,13] = 0.098984
$data[6,14] = 0.001792907809228914
$data[6,15] = 0.001792907809228914
$data[6,16] = 0.05060137967733334
$data[6,17] = 0.4554124170960001
$data[6,18] = 0.0003794174595679246
$data[6,19] = 0.0003794174595679245

$data[7,0] = "FAPs"
$data[7,1] = "Wnt4"
$data[7,2] = "Fzd2"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 1.533623
$data[7,7] = 4.600869
$data[7,8] = 0.2116212878402838
$data[7,9] = 0.2116212878402837
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 2.583925666666666
$data[7,13] = 7.751776999999999
$data[7,14] = 0.140408768272661
$data[7,15] = 0.1404087682726611
$data[7,16] = 3.962767832690333
$data[7,17] = 35.664910494213
$data[7,18] = 0.0297134843659285
$data[7,19] = 0.0297134843659285

$data[8,0] = "M2"
$data[8,1] = "Wnt4"
$data[8,2] = "Fzd2"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 1.389186666666667
$data[8,7] = 4.16756
$data[8,8] = 0.1916908336993843
$data[8,9] = 0.1916908336993843
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 0.08106100000000001
$data[8,13] = 0.243183
$data[8,14] = 0.004404799763312406
$data[8,15] = 0.004404799763312406
$data[8,16] = 0.1126088603866667
$data[8,17] = 1.01347974348
$data[8,18] = 0.0008443597389082058
$data[8,19] = 0.0008443597389082058

$data[9,0] = "M2"
$data[9,1] = "Wnt4"
$data[9,2] = "Fzd2"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 1.389186666666667
$data[9,7] = 4.16756
$data[9,8] = 0.1916908336993843
$data[9,9] = 0.1916908336993843
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 15.70489833333333
$data[9,13] = 47.114695
$data[9,14] = 0.8533935241547975
$data[9,15] = 0.8533935241547975
$data[9,16] = 21.81703536602222
$data[9,17] = 196.3533182942
$data[9,18] = 0.1635877161188888
$data[9,19] = 0.1635877161188888

$data[10,0] = "M2"
$data[10,1] = "Wnt4"
$data[10,2] = "Fzd2"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 1.389186666666667
$data[10,7] = 4.16756
$data[10,8] = 0.1916908336993843
$data[10,9] = 0.1916908336993843
$data[10,10] = 2
$data[10,11] = 0.6666666666666666
$data[10,12] = 0.03299466666666667
$data[10,13] = 0.098984
$data[10,14] = 0.001792907809228914
$data[10,15] = 0.001792907809228914
$data[10,16] = 0.04583575100444444
$data[10,17] = 0.41252175904
$data[10,18] = 0.0003436839926972273
$data[10,19] = 0.0003436839926972273

$data[11,0] = "M2"
$data[11,1] = "Wnt4"
$data[11,2] = "Fzd2"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 1.389186666666667
$data[11,7] = 4.16756
$data[11,8] = 0.1916908336993843
$data[11,9] = 0.1916908336993843
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 2.583925666666666
$data[11,13] = 7.751776999999999
$data[11,14] = 0.140408768272661
$data[11,15] = 0.1404087682726611
$data[11,16] = 3.58955508379111
$data[11,17] = 32.30599575412
$data[11,18] = 0.02691507384889006
$data[11,19] = 0.02691507384889007

$data[12,0] = "sCs"
$data[12,1] = "Wnt4"
$data[12,2] = "Fzd2"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 1.038885
$data[12,7] = 3.116655
$data[12,8] = 0.1433534718884322
$data[12,9] = 0.1433534718884322
$data[12,10] = 2
$data[12,11] = 0.6666666666666666
$data[12,12] = 0.08106100000000001
$data[12,13] = 0.243183
$data[12,14] = 0.004404799763312406
$data[12,15] = 0.004404799763312406
$data[12,16] = 0.08421305698499999
$data[12,17] = 0.757917512865
$data[12,18] = 0.0006314433390441778
$data[12,19] = 0.0006314433390441778

$data[13,0] = "sCs"
$data[13,1] = "Wnt4"
$data[13,2] = "Fzd2"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 1.038885
$data[13,7] = 3.116655
$data[13,8] = 0.1433534718884322
$data[13,9] = 0.1433534718884322
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 15.70489833333333
$data[13,13] = 47.114695
$data[13,14] = 0.8533935241547975
$data[13,15] = 0.8533935241547975
$data[13,16] = 16.315583305025
$data[13,17] = 146.840249745225
$data[13,18] = 0.1223369245746949
$data[13,19] = 0.1223369245746949

$data[14,0] = "sCs"
$data[14,1] = "Wnt4"
$data[14,2] = "Fzd2"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 1.038885
$data[14,7] = 3.116655
$data[14,8] = 0.1433534718884322
$data[14,9] = 0.1433534718884322
$data[14,10] = 2
$data[14,11] = 0.6666666666666666
$data[14,12] = 0.03299466666666667
$data[14,13] = 0.098984
$data[14,14] = 0.001792907809228914
$data[14,15] = 0.001792907809228914
$data[14,16] = 0.03427766427999999
$data[14,17] = 0.30849897852
$data[14,18] = 0.0002570195592288477
$data[14,19] = 0.0002570195592288478

$data[15,0] = "sCs"
$data[15,1] = "Wnt4"
$data[15,2] = "Fzd2"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 1.038885
$data[15,7] = 3.116655
$data[15,8] = 0.1433534718884322
$data[15,9] = 0.1433534718884322
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 2.583925666666666
$data[15,13] = 7.751776999999999
$data[15,14] = 0.140408768272661
$data[15,15] = 0.1404087682726611
$data[15,16] = 2.684401616214999
$data[15,17] = 24.15961454593499
$data[15,18] = 0.02012808441546431
$data[15,19] = 0.02012808441546431

$ws.Range("A2:T17").Value = $data
Write-Host "Updated Wnt4-Fzd2 table with Dr Hou advice values"